$d = $word.ActiveDocument

# --- Part 1: split "ables do not have to use the local" into
#     "ables " + "and functions " + "do not have to use the local"
#     around the existing _GoBack bookmark.
$target = $d.Content.Find.Execute("ables do not have to use the local", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$found = $d.Range(0, 0)
$searchRange = $d.Content
$searchRange.Find.Text = "ables do not have to use the local"
$searchRange.Find.Forward = $true
$searchRange.Find.Wrap = 0
if ($searchRange.Find.Execute()) {
    # searchRange now covers the found text "ables do not have to use the local"
    $start = $searchRange.Start
    $end = $searchRange.End

    # Replace just this run's text with "ables " (first chunk)
    $searchRange.Text = "ables and functions do not have to use the local"
}

Write-Host "done"
